$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.268.69"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").Value = "1.907.78"
$ws.Range("E3").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.78"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5266"
$ws.Range("E7").Value = "  +1.23%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3821"
$ws.Range("E8").Value = "  +1.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07296"
$ws.Range("E9").Value = "  +0.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.11"
$ws.Range("E10").Value = "  +4.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9033"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08187"
$ws.Range("E12").Value = "  -1.64%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.89"
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.363"
$ws.Range("E14").Value = "  +1.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.002"
$ws.Range("E15").Value = "  +0.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008645"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.80"
$ws.Range("E17").Value = "  +1.82%  "

$ws.Range("D18").Value = "1.384.10"
$ws.Range("E18").Value = "  -27.33%  "

$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("D20").Value = "27.307.17"
$ws.Range("E20").Value = "  +0.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.083"

$ws.Range("E22").Value = "  +1.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.529"
$ws.Range("E23").Value = "  +1.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.75"
$ws.Range("E24").Value = "  +2.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.303"
$ws.Range("E25").Value = "  -0.91%  "

$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("E27").Value = "  -0.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.60"
$ws.Range("E28").Value = "  +1.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.838"
$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.836"
$ws.Range("E30").Value = "  -1.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09263"
$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8322"
$ws.Range("E32").Value = "  +4.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05083"
$ws.Range("E33").Value = "  +0.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.232"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.001"
$ws.Range("E35").Value = "  +1.90%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.361"
$ws.Range("E36").Value = "  -1.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.685"
$ws.Range("E37").Value = "  +3.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5829"
$ws.Range("E38").Value = "  +1.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02007"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("E40").Value = "  +0.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.333"
$ws.Range("E41").Value = "  +3.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.553"
$ws.Range("E42").Value = "  -0.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.93"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1525"
$ws.Range("E44").Value = "  +0.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4930"
$ws.Range("E45").Value = "  +1.43%  "

$ws.Range("E46").Value = "  +0.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.16"
$ws.Range("E47").Value = "  +0.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.640"
$ws.Range("E48").Value = "  +0.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.86"
$ws.Range("E49").Value = "  +2.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06211"
$ws.Range("E50").Value = "  +4.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.11"
$ws.Range("E51").Value = "  +0.21%  "
